$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price/volume data to match the latest scrape.
# Cells whose new value is a plain numeric string need the cell
# number format forced to Text ("@") first, otherwise Excel would
# auto-convert the literal into a real number (losing the intended
# text-formatted price representation used throughout this sheet).

# Row 2
$ws.Range("D2").Value = "58.592.95"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.539.19"
$ws.Range("E3").Value = "  +2.02%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.18"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.43"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.87%  "

# Row 9
$ws.Range("D9").Value = "2.538.93"
$ws.Range("E9").Value = "  +1.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0984"
$ws.Range("E10").Value = "  -0.74%  "

# Row 11
$ws.Range("E11").Value = "  -1.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("E12").Value = "  -2.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  -2.47%  "

# Row 14
$ws.Range("D14").Value = "2.978.06"
$ws.Range("E14").Value = "  +1.62%  "

# Row 15
$ws.Range("D15").Value = "58.358.77"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.28"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("D18").Value = "2.530.85"
$ws.Range("E18").Value = "  +1.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.65"
$ws.Range("E20").Value = "  +0.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -0.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  +6.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.03"
$ws.Range("E24").Value = "  +0.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("E25").Value = "  -1.23%  "

# Row 26
$ws.Range("E26").Value = "  +0.30%  "

# Row 27
$ws.Range("E27").Value = "  -0.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  -0.45%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0754"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("E30").Value = "  +1.37%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.42"
$ws.Range("E31").Value = "  -0.87%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  +0.20%  "

# Row 33
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.33"
$ws.Range("E33").Value = "  -0.52%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("E35").Value = "  +0.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.30"
$ws.Range("E36").Value = "  +0.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("E37").Value = "  -5.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.94"
$ws.Range("E38").Value = "  -2.45%  "

# Row 39
$ws.Range("E39").Value = "  +0.69%  "

# Row 40
$ws.Range("E40").Value = "  -0.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.777"
$ws.Range("E41").Value = "  -3.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "279.00"
$ws.Range("E42").Value = "  +0.64%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("E44").Value = "  -3.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("E45").Value = "  +0.61%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.02"
$ws.Range("E46").Value = "  +4.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0918"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0501"
$ws.Range("E48").Value = "  +1.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.79"
$ws.Range("E49").Value = "  -0.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0215"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.13"
$ws.Range("E51").Value = "  -0.59%  "
